$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:C43 with new computed values
$ws.Range("B2").Value2 = 8.655440747545425
$ws.Range("C2").Value2 = 0.707900011048036
$ws.Range("B3").Value2 = 8.774132905390809
$ws.Range("C3").Value2 = 1.404152533376853
$ws.Range("B4").Value2 = 9.161640234321814
$ws.Range("C4").Value2 = 2.05969221672366
$ws.Range("B5").Value2 = 9.220165123224707
$ws.Range("C5").Value2 = 2.756058946153658
$ws.Range("B6").Value2 = 11.37189465995256
$ws.Range("C6").Value2 = 3.775199062348046
$ws.Range("B7").Value2 = 11.48692781600726
$ws.Range("C7").Value2 = 4.574250469192681
$ws.Range("B8").Value2 = 14.83339917805259
$ws.Range("C8").Value2 = 5.306038798351794
$ws.Range("B9").Value2 = 15.11230411077724
$ws.Range("C9").Value2 = 5.998675075632562
$ws.Range("B10").Value2 = 16.90935188785546
$ws.Range("C10").Value2 = 6.829170393289028
$ws.Range("B11").Value2 = 17.00663096371362
$ws.Range("C11").Value2 = 7.630719736501384
$ws.Range("B12").Value2 = 20.88084421830928
$ws.Range("C12").Value2 = 8.494054961347006
$ws.Range("B13").Value2 = 21.33981148266859
$ws.Range("C13").Value2 = 9.146968971703915
$ws.Range("B14").Value2 = 21.41912674689676
$ws.Range("C14").Value2 = 9.933491801155858
$ws.Range("B15").Value2 = 24.96119769513005
$ws.Range("C15").Value2 = 10.72489668234038
$ws.Range("B16").Value2 = 25.02477168113873
$ws.Range("C16").Value2 = 11.49069797649732
$ws.Range("B17").Value2 = 25.38148241481064
$ws.Range("C17").Value2 = 12.70703163302527
$ws.Range("B18").Value2 = 25.80136852157638
$ws.Range("C18").Value2 = 13.40721855107083
$ws.Range("B19").Value2 = 26.62430429673192
$ws.Range("C19").Value2 = 14.33904275983722
$ws.Range("B20").Value2 = 31.28365183642197
$ws.Range("C20").Value2 = 15.13756516759278
$ws.Range("B21").Value2 = 31.73580809556347
$ws.Range("C21").Value2 = 16.02461502984101
$ws.Range("B22").Value2 = 31.7883100958229
$ws.Range("C22").Value2 = 16.69584912948651
$ws.Range("B23").Value2 = 32.83302442758634
$ws.Range("C23").Value2 = 17.29244698511091
$ws.Range("B24").Value2 = 32.93053329523873
$ws.Range("C24").Value2 = 18.01967643676003
$ws.Range("B25").Value2 = 33.85998950212554
$ws.Range("C25").Value2 = 18.9154698499765
$ws.Range("B26").Value2 = 35.33772968032522
$ws.Range("C26").Value2 = 19.78680256473768
$ws.Range("B27").Value2 = 38.88832562868661
$ws.Range("C27").Value2 = 20.57943623161137
$ws.Range("B28").Value2 = 38.92632501388512
$ws.Range("C28").Value2 = 21.24935684191707
$ws.Range("B29").Value2 = 44.59948247572704
$ws.Range("C29").Value2 = 21.98693187036435
$ws.Range("B30").Value2 = 44.68661802200463
$ws.Range("C30").Value2 = 22.73611222050961
$ws.Range("B31").Value2 = 45.77237404158237
$ws.Range("C31").Value2 = 23.54550444671809
$ws.Range("B32").Value2 = 59.75624638292181
$ws.Range("C32").Value2 = 24.30493243624291
$ws.Range("B33").Value2 = 59.81130430988183
$ws.Range("C33").Value2 = 25.19035318066064
$ws.Range("B34").Value2 = 67.03414541096291
$ws.Range("C34").Value2 = 25.82903171167909
$ws.Range("B35").Value2 = 67.63448670982559
$ws.Range("C35").Value2 = 26.70390099532375
$ws.Range("B36").Value2 = 71.18229273445434
$ws.Range("C36").Value2 = 27.36050330513353
$ws.Range("B37").Value2 = 71.2408045326759
$ws.Range("C37").Value2 = 28.14263578384396
$ws.Range("B38").Value2 = 88.2132499314002
$ws.Range("C38").Value2 = 28.88185357846553
$ws.Range("B39").Value2 = 90.48472526696902
$ws.Range("C39").Value2 = 29.80205304393597
$ws.Range("B40").Value2 = 90.7240637115648
$ws.Range("C40").Value2 = 30.56439846491084
$ws.Range("B41").Value2 = 93.83070237118476
$ws.Range("C41").Value2 = 31.38826187116699
$ws.Range("B42").Value2 = 94.127209970392
$ws.Range("C42").Value2 = 32.56127090299168
$ws.Range("B43").Value2 = 94.34898901636633
$ws.Range("C43").Value2 = 33.33378394548637

# Remove rows 44-48 (former A=42..46 data), which no longer exist
$ws.Range("A44:A48").EntireRow.Delete() | Out-Null

